$wb = $excel.ActiveWorkbook

# sigma_010 sheet: update Rows index (0-based -> 1-based) and refreshed Noisy/NLM-LBP values
$ws = $wb.Worksheets.Item("sigma_010")
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = 27.74812892313493
$ws.Cells.Item(2,3).Value = 29.66486156874257
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = 27.75698092535827
$ws.Cells.Item(3,3).Value = 29.65080522938081
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = 27.77940329728607
$ws.Cells.Item(4,3).Value = 29.67237505860282
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = 27.76010084592922
$ws.Cells.Item(5,3).Value = 29.67068022056172
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = 27.78125461783738
$ws.Cells.Item(6,3).Value = 29.67844141706296
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = 27.76238316518699
$ws.Cells.Item(7,3).Value = 29.68382041937317
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = 27.77356551444721
$ws.Cells.Item(8,3).Value = 29.64570177656125
$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = 27.78067889727183
$ws.Cells.Item(9,3).Value = 29.67942536583637
$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = 27.7931705289115
$ws.Cells.Item(10,3).Value = 29.69683099216793
$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = 27.79610936803971
$ws.Cells.Item(11,3).Value = 29.68008838992794
$ws.Cells.Item(12,2).Value = 27.77317760834031
$ws.Cells.Item(12,3).Value = 29.67230304382175

# sigma_025 sheet: update Rows index (0-based -> 1-based) and refreshed Noisy/NLM-LBP values
$ws = $wb.Worksheets.Item("sigma_025")
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = 19.71749016998847
$ws.Cells.Item(2,3).Value = 26.44959565633147
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = 19.73032195981172
$ws.Cells.Item(3,3).Value = 26.42777620902709
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = 19.70291114544299
$ws.Cells.Item(4,3).Value = 26.41437967370847
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = 19.71614988475168
$ws.Cells.Item(5,3).Value = 26.41804150610024
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = 19.70993381839491
$ws.Cells.Item(6,3).Value = 26.43177649193372
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = 19.71819884433467
$ws.Cells.Item(7,3).Value = 26.44762132471604
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = 19.71493673140955
$ws.Cells.Item(8,3).Value = 26.42083312970303
$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = 19.7038778279932
$ws.Cells.Item(9,3).Value = 26.42111397386038
$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = 19.73886991892218
$ws.Cells.Item(10,3).Value = 26.44026110164843
$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = 19.72285312578322
$ws.Cells.Item(11,3).Value = 26.45735935492105
$ws.Cells.Item(12,2).Value = 19.71755434268326
$ws.Cells.Item(12,3).Value = 26.43287584219499

# sigma_050 sheet: update Rows index (0-based -> 1-based) and refreshed Noisy/NLM-LBP values
$ws = $wb.Worksheets.Item("sigma_050")
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = 14.66408655001434
$ws.Cells.Item(2,3).Value = 21.16540901400988
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = 14.65992773565036
$ws.Cells.Item(3,3).Value = 21.18259530504783
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = 14.66482750758352
$ws.Cells.Item(4,3).Value = 21.21897729071021
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = 14.66651342633821
$ws.Cells.Item(5,3).Value = 21.12655249795798
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = 14.66991808155307
$ws.Cells.Item(6,3).Value = 21.20434807015733
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = 14.67048354894018
$ws.Cells.Item(7,3).Value = 21.20663313012475
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = 14.68349634033455
$ws.Cells.Item(8,3).Value = 21.20498154390098
$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = 14.68181423843806
$ws.Cells.Item(9,3).Value = 21.13653338148296
$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = 14.65812643244132
$ws.Cells.Item(10,3).Value = 21.19586775471614
$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = 14.67150906180261
$ws.Cells.Item(11,3).Value = 21.18548318772741
$ws.Cells.Item(12,2).Value = 14.66907029230962
$ws.Cells.Item(12,3).Value = 21.18273811758355
